# "Fixed recommendations score & User likes"
#
# The 9 users keep their existing rows (2-10); only their
# Video games / Sport / Music percentages are corrected to a clean
# single-category preference (100/0/0, cycled across the three columns).
# The "Comment" column (header + explanatory notes) is no longer used and
# is cleared out, and the worksheet selection moves to B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected preference scores (Video games %, Sport %, Music %) -------
$scores = @{
    2  = @(100, 0, 0)     # Govinda Dimosthenis
    3  = @(0, 100, 0)     # Azarel Feodosiy
    4  = @(0, 0, 100)     # Jarah Aaron
    5  = @(100, 0, 0)     # Peder Oddmund
    6  = @(0, 100, 0)     # Ashok Walganus
    7  = @(0, 0, 100)     # Royle Faust
    8  = @(100, 0, 0)     # Azhar Drake
    9  = @(0, 100, 0)     # Artur Emil
    10 = @(0, 0, 100)     # Yoshirou Harvie
}

foreach ($row in $scores.Keys) {
    $vals = $scores[$row]
    $ws.Cells.Item($row, 2).Value2 = $vals[0]
    $ws.Cells.Item($row, 3).Value2 = $vals[1]
    $ws.Cells.Item($row, 4).Value2 = $vals[2]
}

# --- Drop the now-unused "Comment" column text ----------------------------
# (Value2 = "" is used instead of ClearContents() because E2:E4/E6:E7/E8:E9
# are merged ranges and ClearContents() does not propagate through the
# merge on this host - direct value assignment on the anchor cell does.)
$ws.Range("E1").Value2 = ""
$ws.Range("E2").Value2 = ""
$ws.Range("E3").Value2 = ""
$ws.Range("E4").Value2 = ""
$ws.Range("E5").Value2 = ""
$ws.Range("E6").Value2 = ""
$ws.Range("E8").Value2 = ""
$ws.Range("E10").Value2 = ""

# --- Move the active selection ---------------------------------------------
$ws.Range("B11").Select()
